$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Boolean": expand the two combined-vehicle-type CSV rows into six
# per-vehicle-type rows each, and clear/select a couple cells the way the
# authored workbook ended up (adds six trailing blank rows too).
# ---------------------------------------------------------------------------
$wsBool = $wb.Worksheets.Item("Boolean")

# Row 17 currently holds "trans/BVTQaZ/BVTQaZ.csv" - expand to 6 rows.
$wsBool.Rows.Item(17).Resize(5).Insert()
$wsBool.Cells.Item(17, 1).Value = "trans/BVTQaZ/BVTQaZ-LDVs.csv"
$wsBool.Cells.Item(18, 1).Value = "trans/BVTQaZ/BVTQaZ-HDVs.csv"
$wsBool.Cells.Item(19, 1).Value = "trans/BVTQaZ/BVTQaZ-aircraft.csv"
$wsBool.Cells.Item(20, 1).Value = "trans/BVTQaZ/BVTQaZ-rail.csv"
$wsBool.Cells.Item(21, 1).Value = "trans/BVTQaZ/BVTQaZ-ships.csv"
$wsBool.Cells.Item(22, 1).Value = "trans/BVTQaZ/BVTQaZ-motorbikes.csv"

# Row 26 (after the insert above) now holds "trans/VTQaZ/VTQaZ.csv" - expand
# to 6 rows too.
$wsBool.Rows.Item(26).Resize(5).Insert()
$wsBool.Cells.Item(26, 1).Value = "trans/VTQaZ/VTQaZ-LDVs.csv"
$wsBool.Cells.Item(27, 1).Value = "trans/VTQaZ/VTQaZ-HDVs.csv"
$wsBool.Cells.Item(28, 1).Value = "trans/VTQaZ/VTQaZ-aircraft.csv"
$wsBool.Cells.Item(29, 1).Value = "trans/VTQaZ/VTQaZ-rail.csv"
$wsBool.Cells.Item(30, 1).Value = "trans/VTQaZ/VTQaZ-ships.csv"
$wsBool.Cells.Item(31, 1).Value = "trans/VTQaZ/VTQaZ-motorbikes.csv"

# Six trailing blank rows (33-38) were added after the last data row (32).
$wsBool.Rows.Item(33).Resize(6).Insert()

# Scroll/selection state ends up parked near the bottom of the new list.
$wsBool.Range("A32").Select()

# ---------------------------------------------------------------------------
# Sheet "About": the closing note now points at the new "data types" list
# instead of the old fuels CSV path.
# ---------------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A7").Value = "InputData pathnames of CSV files with values constrained to specific data types"

# ---------------------------------------------------------------------------
# Sheet "Integer": selection parked at A13.
# ---------------------------------------------------------------------------
$wsInt = $wb.Worksheets.Item("Integer")
$wsInt.Range("A13").Select()

# ---------------------------------------------------------------------------
# Make "About" the active tab (it picks up tabSelected="1" on save).
# ---------------------------------------------------------------------------
$wsAbout.Activate()
$wsAbout.Range("A1").Select()
